$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44894
$ws.Range("L2").Value = 'Segunda'
$ws.Range("M2").Value = 130
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19462
$ws.Range("Q2").Value = '$/caja 16 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1216
$ws.Range("T2").Value = 16

$ws.Range("D3").Value = 44533
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 140
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("Q3").Value = '$/caja 10 kilos'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 1450
$ws.Range("T3").Value = 10

$ws.Range("D4").Value = 44524
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 27000
$ws.Range("O4").Value = 28000
$ws.Range("P4").Value = 27500
$ws.Range("Q4").Value = '$/bandeja 18 kilos'
$ws.Range("R4").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S4").Value = 1528
$ws.Range("T4").Value = 18

$ws.Range("D6").Value = 44881
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 41000
$ws.Range("O6").Value = 42000
$ws.Range("P6").Value = 41500
$ws.Range("Q6").Value = '$/bandeja 18 kilos'
$ws.Range("R6").Value = 'Región de Coquimbo'
$ws.Range("S6").Value = 2306
$ws.Range("T6").Value = 18

$ws.Range("D7").Value = 44160
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 24000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 24500
$ws.Range("Q7").Value = '$/bandeja 18 kilos'
$ws.Range("R7").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S7").Value = 1361
$ws.Range("T7").Value = 18

$ws.Range("D8").Value = 44169
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 250
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 21000
$ws.Range("Q8").Value = '$/bandeja 18 kilos'
$ws.Range("R8").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S8").Value = 1167
$ws.Range("T8").Value = 18

$ws.Range("D9").Value = 44895
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 130
$ws.Range("N9").Value = 19000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 19462
$ws.Range("Q9").Value = '$/caja 16 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1216
$ws.Range("T9").Value = 16

$ws.Range("D10").Value = 44917
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 250
$ws.Range("N10").Value = 20000
$ws.Range("O10").Value = 23000
$ws.Range("P10").Value = 21800
$ws.Range("Q10").Value = '$/caja 18 kilos'
$ws.Range("R10").Value = 'Región de Coquimbo'
$ws.Range("S10").Value = 1211
$ws.Range("T10").Value = 18

$ws.Range("D11").Value = 44544
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 21000
$ws.Range("Q11").Value = '$/bandeja 18 kilos'
$ws.Range("R11").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S11").Value = 1167
$ws.Range("T11").Value = 18

$ws.Range("D12").Value = 44901
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 17500
$ws.Range("Q12").Value = '$/bandeja 18 kilos'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 972
$ws.Range("T12").Value = 18

$ws.Range("D13").Value = 44545
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 25000
$ws.Range("P13").Value = 24500
$ws.Range("Q13").Value = '$/bandeja 18 kilos'
$ws.Range("R13").Value = 'Región de Coquimbo'
$ws.Range("S13").Value = 1361
$ws.Range("T13").Value = 18

$ws.Range("D14").Value = 44880
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 33000
$ws.Range("O14").Value = 34000
$ws.Range("P14").Value = 33500
$ws.Range("Q14").Value = '$/caja 10 kilos'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 3350
$ws.Range("T14").Value = 10

